$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 4442.76792576475
$ws.Range("C2").Value = 5521.13092777069
$ws.Range("D2").Value = 5780
$ws.Range("F2").Value = 72.1570401669144

# Row 3
$ws.Range("B3").Value = 4416.29615663236
$ws.Range("C3").Value = 5481.90082850347
$ws.Range("F3").Value = 150.532189952963

# Row 4
$ws.Range("B4").Value = 4361.04011919423
$ws.Range("C4").Value = 5296.5275635083
$ws.Range("F4").Value = 139.885101054753

# Row 5
$ws.Range("B5").Value = 4320.7661481526
$ws.Range("C5").Value = 5108.30143503242
$ws.Range("F5").Value = 130.390191161659

# Row 6
$ws.Range("B6").Value = 5308.36389769672
$ws.Range("C6").Value = 4738.81231410109
$ws.Range("F6").Value = 143.822644225182

# Row 7
$ws.Range("B7").Value = 1365.33360290578
$ws.Range("C7").Value = 3203.79339216272
$ws.Range("F7").Value = 87.3475197190394

# Row 8
$ws.Range("C8").Value = 3402.34994971906
$ws.Range("F8").Value = 89.4570708123993

# Row 9
$ws.Range("C9").Value = 5598.09677609297
$ws.Range("F9").Value = 212.002472133276

# Row 10
$ws.Range("C10").Value = 5605.17118803607
$ws.Range("F10").Value = 212.297239297572

# Row 11
$ws.Range("C11").Value = 5783.50423013124
$ws.Range("F11").Value = 220.542405926537

# Row 12
$ws.Range("C12").Value = 5991.50209463559
$ws.Range("F12").Value = 229.208983614218

# Row 13
$ws.Range("C13").Value = 5407.86951544969
$ws.Range("F13").Value = 204.890959481472

# Row 14
$ws.Range("C14").Value = 3852.20362797654
$ws.Range("F14").Value = 123.482007503138

# Row 15
$ws.Range("C15").Value = 3871.48332234216
$ws.Range("F15").Value = 123.875992682245

$wb.Save()
